$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection (active cell) to G7
$ws.Range("G7").Select()

# D3: previously text "M0012", now becomes numeric value 1
$ws.Range("D3").Value = 1

# G3: previously empty (just formatted as text), now gets phone number text
$ws.Range("G3").Value = "082345675423"

# O3 stays "N" (its underlying shared string index shifts automatically
# since "M0012" is no longer referenced/used in the sheet)
$ws.Range("O3").Value = "N"
